$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("P4").Value = "320018501311"
$ws.Range("Q4").Value = "$76.67"
$ws.Range("R4").Value = "FAIL"
